$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$rows = @(19, 22, 25)

foreach ($r in $rows) {
    $hCell = $ws.Range("H" + $r)
    $jCell = $ws.Range("J" + $r)

    # Grab the text value currently stored in H
    $val = $hCell.Value2

    # Copy H's formatting onto J (clipboard copy + paste-special formats only,
    # so the existing merge areas of H and J are left untouched)
    $hCell.Copy()
    $jCell.PasteSpecial($xlPasteFormats)

    # Move the value across: J gets H's old text, H becomes blank
    $jCell.Value2 = $val
    $hCell.Value2 = ""
}

$excel.CutCopyMode = 0
